$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from an existing header cell (e.g. E1) to the new header cells
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

$ws.Range("E1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)

$ws.Range("F2").Value = $false
$ws.Range("G2").Value = $false
$ws.Range("H2").Value = $false
